$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "58.103.92"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -4.14%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.977.94"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.10%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "559.99"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.06%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "133.80"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.91%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("E8").Value = "  +3.00%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.971.04"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.21%  "
$ws.Range("E10").Value = "  -2.93%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "4.87"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -5.55%  "
$ws.Range("E12").Value = "  +1.86%  "
$ws.Range("E13").Value = "  +0.08%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "33.06"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.92%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.468.23"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.18%  "
$ws.Range("E17").Value = "  +7.35%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.972.28"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.41%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "58.040.27"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -4.06%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "421.74"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.60%  "
$ws.Range("E21").Value = "  +0.83%  "
$ws.Range("E22").Value = "  +3.42%  "
$ws.Range("E23").Value = "  -0.40%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.09"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.68%  "
$ws.Range("E25").Value = "  +0.20%  "
$ws.Range("E26").Value = "  -0.03%  "
$ws.Range("E27").Value = "  +0.01%  "
$ws.Range("E28").Value = "  -2.23%  "
$ws.Range("E29").Value = "  +4.00%  "
$ws.Range("E30").Value = "  +5.15%  "
$ws.Range("B31").Value = "EthereumClassic"
$ws.Range("C31").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "25.42"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.11%  "
$ws.Range("B32").Value = "NEARProtocol"
$ws.Range("C32").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.12"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.05%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.100"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +6.70%  "
$ws.Range("E34").Value = "  -1.16%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.67"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.95%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.947"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.54%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0₃0702"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +4.97%  "
$ws.Range("E38").Value = "  -2.77%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.70"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.41%  "
$ws.Range("E40").Value = "  +2.89%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0353"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.40%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.109"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.95%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "379.56"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.18%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.690.36"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.94%  "
$ws.Range("E45").Value = "  -0.05%  "
$ws.Range("E46").Value = "  +2.87%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "121.90"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.33%  "
$ws.Range("E48").Value = "  +2.58%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.01"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.16%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "23.61"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.24%  "
$ws.Range("E51").Value = "  -0.42%  "
